$wb = $excel.ActiveWorkbook

# Sheet: labor_incmon_imp_stochastic_reg
$ws1 = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws1.Range("G3").Value = 2080805.9128611058
$ws1.Range("H3").Value = 299415.71875
$ws1.Range("I3").Value = 538948.3125
$ws1.Range("J3").Value = 1461299
$ws1.Range("K3").Value = 2653803.75
$ws1.Range("L3").Value = 3880000

# Sheet: labor_jubpenimp_stochastic_reg
$ws2 = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws2.Range("G3").Value = 1542192.2846315436

# Sheet: nonlabor_imp_stochastic_reg
$ws3 = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws3.Range("G3").Value = 9105707.0312418621
$ws3.Range("L3").Value = 900764.8125

# Sheet: labor_beneimp_stochastic_reg
$ws4 = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws4.Range("G3").Value = 866425.13373399491
$ws4.Range("H3").Value = 60000
$ws4.Range("I3").Value = 199610.484375
$ws4.Range("J3").Value = 349318.34375
$ws4.Range("K3").Value = 945916.4375
$ws4.Range("L3").Value = 2000000
